$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# Step 1: Insert a new bordered paragraph right before the current paragraph 2
# ("Remember that this is a public repository ...") containing the new
# "Github Username: LitheDev" bold text (with spell-check proofErr markers,
# matching a Word autocorrect/spell-check artifact).
# ---------------------------------------------------------------------------
$target = $d.Paragraphs(2).Range
$target.InsertParagraphBefore()
$ghPara = $d.Paragraphs(2).Range

$ghXml = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?>' +
  '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">' +
  '<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData>' +
  '<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">' +
  '<w:body><w:p>' +
  '<w:pPr><w:pBdr>' +
  '<w:top w:val="single" w:sz="4" w:space="1" w:color="auto"/>' +
  '<w:left w:val="single" w:sz="4" w:space="4" w:color="auto"/>' +
  '<w:bottom w:val="single" w:sz="4" w:space="1" w:color="auto"/>' +
  '<w:right w:val="single" w:sz="4" w:space="4" w:color="auto"/>' +
  '</w:pBdr><w:rPr><w:b/><w:bCs/></w:rPr></w:pPr>' +
  '<w:proofErr w:type="spellStart"/>' +
  '<w:r><w:rPr><w:b/><w:bCs/></w:rPr><w:t>Github</w:t></w:r>' +
  '<w:proofErr w:type="spellEnd"/>' +
  '<w:r><w:rPr><w:b/><w:bCs/></w:rPr><w:t xml:space="preserve"> Username: </w:t></w:r>' +
  '<w:proofErr w:type="spellStart"/>' +
  '<w:r><w:rPr><w:b/><w:bCs/></w:rPr><w:t>LitheDev</w:t></w:r>' +
  '<w:proofErr w:type="spellEnd"/>' +
  '</w:p></w:body></w:document>' +
  '</pkg:xmlData></pkg:part></pkg:package>'

$ghPara.InsertXML($ghXml)

# ---------------------------------------------------------------------------
# Step 2: Insert a new bordered paragraph before the current paragraph 6
# ("Once you've changed follow the next step ...") containing the text that
# used to live in the "Add some comments about Version management" paragraph.
# ---------------------------------------------------------------------------
$beforeOnce = $d.Paragraphs(6).Range
$beforeOnce.InsertParagraphBefore()
$addPara = $d.Paragraphs(6).Range

$addXml = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?>' +
  '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">' +
  '<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData>' +
  '<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">' +
  '<w:body><w:p>' +
  '<w:pPr><w:pBdr>' +
  '<w:top w:val="single" w:sz="4" w:space="1" w:color="auto"/>' +
  '<w:left w:val="single" w:sz="4" w:space="4" w:color="auto"/>' +
  '<w:bottom w:val="single" w:sz="4" w:space="1" w:color="auto"/>' +
  '<w:right w:val="single" w:sz="4" w:space="4" w:color="auto"/>' +
  '</w:pBdr></w:pPr>' +
  '<w:r><w:t xml:space="preserve">Add some comments about Version management </w:t></w:r>' +
  '<w:r><w:t>outside this border</w:t></w:r>' +
  '<w:r><w:t>, or just add some text so there is a change to this file.</w:t></w:r>' +
  '</w:p></w:body></w:document>' +
  '</pkg:xmlData></pkg:part></pkg:package>'

$addPara.InsertXML($addXml)

# ---------------------------------------------------------------------------
# Step 3: Empty out the paragraph that originally held the "Add some
# comments about Version management ..." text (now paragraph 4), leaving its
# border intact but with no runs.
# ---------------------------------------------------------------------------
$oldAddPara = $d.Paragraphs(4).Range
$clearRange = $d.Range($oldAddPara.Start, $oldAddPara.End - 1)
if ($clearRange.Start -lt $clearRange.End) {
  $clearRange.Delete()
}

# ---------------------------------------------------------------------------
# Step 4: Replace the "…" paragraph (now paragraph 8) with the Version
# Management essay paragraph.
# ---------------------------------------------------------------------------
$ellipsisPara = $d.Paragraphs(8).Range
$essayRange = $d.Range($ellipsisPara.Start, $ellipsisPara.End - 1)
$essayRange.Text = 'Version Management is a broad term that covers a few key elements when creating a system with relation to programming. It involves keeping track of older versions of a system, newer versions of a system and deviations in that system (such as upgrades to newer UI’s, etc). This is done throughout the development lifecycle to keep track of stable versions of software while allowing for experimental branches to be created and tested in isolation. Systems such as git that allow for tracking of changes and version history are beneficial when multiple programmers are working on the same system. Moreover, if two programmers happen to be working on the same code file at the same time, while traditional systems (such as local filing systems in windows) would overwrite the file with the newer version (after prompted). Systems such as git allow for parallel programming and insertion of new code, instead of complete overriding of other team members work. '
